# Apply updated cryptocurrency price/volume data (and two name/link/price swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''23.240.47'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  -0.81%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''1.617.74'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -0.79%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').Value = '''1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '''  +0.70%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''1.003'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +0.56%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''303.32'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  -1.50%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = '''0.3770'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '''  -0.33%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''51.92'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '''  -2.48%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D10').Value = '''1.224'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  -4.79%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('B11').Value = '''BinanceUSD'
$ws.Range('B11').Style = "Normal"
$ws.Range('C11').Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C11').Style = "Normal"
$ws.Range('D11').Value = '''1.005'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''  +0.72%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('B12').Value = '''Dogecoin'
$ws.Range('B12').Style = "Normal"
$ws.Range('C12').Value = '''https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('C12').Style = "Normal"
$ws.Range('D12').Value = '''0.08054'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  -1.79%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''22.50'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''  -3.71%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = '''6.522'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''  -2.42%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''0.00001240'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -2.16%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''7.187'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  -3.90%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''1.621.77'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  -0.46%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''93.22'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  -1.85%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''0.06913'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  -0.37%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''17.85'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  -3.22%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''1.003'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  +0.61%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''6.419'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -2.78%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''23.260.55'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  -0.77%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''12.69'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  -2.39%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''3.188'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +1.32%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''2.446'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  +0.27%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''21.09'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '''  -1.61%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''148.82'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  -1.65%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = '''5.280'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  -0.69%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = '''134.49'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  -1.74%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''2.295'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  -5.71%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''1.803.66'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  -0.21%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = '''6.716'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '''  -3.28%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''10.80'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  +2.96%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''0.9420'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  -4.22%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = '''0.02800'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '''  -0.37%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = '''0.2517'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '''  -0.94%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''0.08823'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  -0.36%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = '''6.067'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '''  -3.50%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''0.07079'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  -5.47%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = '''  -3.36%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''0.7003'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  -2.45%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = '''16.19'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  -0.18%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''12.17'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  -5.21%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('B45').Value = '''Frax'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').Value = '''1.002'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  +0.54%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('B46').Value = '''Decentraland'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').Value = '''https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').Value = '''0.6413'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  -3.40%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''2.304'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  -2.74%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = '''3.980'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '''  -1.53%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = '''0.07964'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '''  -0.92%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = '''1.196'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '''  -1.76%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''125.18'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  -5.93%  '
$ws.Range('E51').Style = "Normal"
